$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "CJ20" to "OUVP"
$ws.Name = "OUVP"

# Move the active cell/selection from J26 to H20
$ws.Range("H20").Select()
